# "Personal und LPE geändert in Datenbasis"
# Update the LPE/Personal figures on the "Tabelle1" worksheet of the
# Datenbasis workbook, and move the current selection/view as the author
# left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Update the changed figures ---------------------------------------
$ws.Range("J10").Value = 0.69   # was 0.79
$ws.Range("G11").Value = 100    # was 150
$ws.Range("J11").Value = 0.84   # was 0.89
$ws.Range("G12").Value = 100    # was 200

# --- Update the saved view/selection state ------------------------------
$ws.Activate()
$ws.Range("G13").Select()
$excel.ActiveWindow.ScrollColumn = 4   # scroll so column D is leftmost (topLeftCell D1)
$excel.ActiveWindow.ScrollRow = 1
